$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The raw data (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio
# ponderado, Precio $/Kg) for rows 3, 4 and 5 gets cyclically shifted down by
# one row (row3 -> row4, row4 -> row5, row5 -> row3), since every other
# column is identical across these three rows.

$ws.Range("D3").Value = 44257
$ws.Range("M3").Value = 100

$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

$ws.Range("D5").Value = 44252
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("S5").Value = 750
